# Insert a new row at position 89, shifting existing rows 89-136 down to 90-137,
# and populate the new row 89 with the new weekly price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("89:89").Insert()

$row = $ws.Range("A89:T89")
$row.Cells.Item(1, 1).Value = 4
$row.Cells.Item(1, 2).Value = "Feria Lagunitas de Puerto Montt"
$row.Cells.Item(1, 3).Value = "Los Lagos"
$row.Cells.Item(1, 4).Value = 44460
$row.Cells.Item(1, 5).Value = 10
$row.Cells.Item(1, 6).Value = "Fruta"
$row.Cells.Item(1, 7).Value = 100108
$row.Cells.Item(1, 8).Value = "Tropicales y subtropicales"
$row.Cells.Item(1, 9).Value = 100108005
$row.Cells.Item(1, 10).Value = "Piña"
$row.Cells.Item(1, 11).Value = "Caramelo"
$row.Cells.Item(1, 12).Value = "Primera"
$row.Cells.Item(1, 13).Value = 60
$row.Cells.Item(1, 14).Value = 22000
$row.Cells.Item(1, 15).Value = 22000
$row.Cells.Item(1, 16).Value = 22000
$row.Cells.Item(1, 17).Value = "$/caja 12 unidades"
$row.Cells.Item(1, 18).Value = "Ecuador"
$row.Cells.Item(1, 19).Value = 1833
$row.Cells.Item(1, 20).Value = 12
